$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 367, shifting existing rows (367..502) down to (368..503).
$ws.Rows.Item(367).Insert()

# The newly inserted row 367 is a duplicate of the (now shifted-down) former row 367,
# which now lives at row 368 - copy its contents into the new row 367, then
# overwrite the Fecha (D) and Volumen (J) values for the new record.
$srcRow = $ws.Rows.Item(368)
$dstRow = $ws.Rows.Item(367)
$srcRow.Copy()
$dstRow.PasteSpecial(-4104)

$ws.Cells.Item(367, 4).Value = 45009
$ws.Cells.Item(367, 10).Value = 5000
